$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the simulation-results data (rows 2-11, columns B/C/D) ---
# Column A (threshold values) is unchanged.
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 48.02

$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 48.06

$ws.Range("B4").Value = 99.99
$ws.Range("C4").Value = 0.01
$ws.Range("D4").Value = 47.97

$ws.Range("B5").Value = 96.93
$ws.Range("C5").Value = 3.07
$ws.Range("D5").Value = 47.99

$ws.Range("B6").Value = 37.68
$ws.Range("C6").Value = 62.32
$ws.Range("D6").Value = 48.04

$ws.Range("B7").Value = 0.66
$ws.Range("C7").Value = 99.34
$ws.Range("D7").Value = 48.03

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 100
$ws.Range("D8").Value = 47.99

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 100
$ws.Range("D9").Value = 47.98

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 100
$ws.Range("D10").Value = 48.05

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 100
$ws.Range("D11").Value = 47.98

# --- Strip the big bespoke header/cell formatting (big fonts, borders,
#     26pt row height) that the re-exported workbook no longer carries;
#     cells fall back to the default "Normal" style/row height. ---
$ws.Rows("1:11").ClearFormats()
$ws.Columns("A:D").AutoFit()
$ws.Rows("1:11").AutoFit()

# --- Remove the stale hidden "_xlchart.*" defined names left over from
#     the charting add-in; the refreshed workbook no longer references them. ---
$names = @()
foreach ($n in $wb.Names) {
    $names += $n
}
foreach ($n in $names) {
    $n.Delete()
}

# --- Reposition/resize the chart to track the new (narrower) column
#     layout and selection state the author left the sheet in. ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 326.875
$co.Top = 3
$co.Width = 1064.1875
$co.Height = 647

[void]$ws.Range("V15:W15").Select()
$excel.ActiveWindow.Zoom = 100
